$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header label from "Price 2021-05-23" to "Price 2021-05-30"
$ws.Range("G1").Value = "Price 2021-05-30"

# Update price values
$ws.Range("G4").Value = 19.99
$ws.Range("G5").Value = 19.99
$ws.Range("G8").Value = 292.3
